$d = $word.ActiveDocument

$replacements = @(
    @{old = "37×21=777"; new = "57×14=798"},
    @{old = "41×27=1107"; new = "34×37=1258"},
    @{old = "68×93=6324"; new = "26×33=858"},
    @{old = "29×13=377"; new = "51×88=4488"},
    @{old = "62×99=6138"; new = "63×71=4473"},
    @{old = "23×51=1173"; new = "11×68=748"},
    @{old = "66×85=5610"; new = "31×21=651"},
    @{old = "41×29=1189"; new = "73×24=1752"},
    @{old = "33×65=2145"; new = "33×66=2178"},
    @{old = "52×17=884"; new = "92×74=6808"},
    @{old = "47×36=1692"; new = "35×57=1995"},
    @{old = "68×51=3468"; new = "43×37=1591"},
    @{old = "61×33=2013"; new = "79×64=5056"},
    @{old = "24×99=2376"; new = "87×46=4002"},
    @{old = "51×78=3978"; new = "20×90=1800"},
    @{old = "40×64=2560"; new = "35×96=3360"},
    @{old = "44×62=2728"; new = "73×13=949"},
    @{old = "17×42=714"; new = "92×50=4600"},
    @{old = "43×99=4257"; new = "38×29=1102"},
    @{old = "64×47=3008"; new = "19×76=1444"},
    @{old = "18×68=1224"; new = "49×71=3479"},
    @{old = "69×32=2208"; new = "20×45=900"},
    @{old = "41×20=820"; new = "27×12=324"},
    @{old = "18×35=630"; new = "49×91=4459"},
    @{old = "76×35=2660"; new = "84×27=2268"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
